$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move 1: "Leesub Sirln" (currently row 142) moves up to right after
#     "Bail Organa" / right before "Jenny" (row 108). Implemented as
#     Copy -> Insert at destination -> Delete the now-duplicated source row
#     (the classic COM "move row" idiom, since Cut+Insert isn't reliable here).
$ws.Rows(142).Copy()
$ws.Rows(108).Insert()
$ws.Rows(143).Delete()

# --- Move 2: "Momaw Nadon" (originally row 137, now shifted to row 138 by
#     the insert above) moves down to right after "Trech Molock" (row 141),
#     i.e. right before "Beedo" (row 143).
$ws.Rows(138).Copy()
$ws.Rows(143).Insert()
$ws.Rows(138).Delete()
